# The workbook contains two adjacent species-observation records (rows 43
# and 44) that were swapped: everything that used to be row 44 is now row
# 43, and everything that used to be row 43 is now row 44. Implement this
# by swapping the cell contents of the two rows, column by column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding genuine numbers (must stay numeric, not become text).
$numericCols = @("A", "B", "E", "Q", "R", "S")
# Columns holding booleans.
$booleanCols = @("AD", "AE", "AG")
# All other populated columns hold text (some of which look like dates,
# e.g. "2026-01-20", and must be kept as literal text rather than being
# reinterpreted as a date serial by Excel on assignment).
$textCols = @("D", "F", "G", "H", "I", "J", "P", "T", "U", "V", "W", "Y", "Z", "AA", "AB", "AC", "AJ", "AK", "AO", "AT", "AW", "AX", "AY")

$allCols = $numericCols + $booleanCols + $textCols

function Set-CellFrom($destRange, $value, $isText) {
    if ($value -eq $null) {
        $destRange.ClearContents()
        return
    }
    if ($isText) {
        # Force text interpretation so date-looking strings ("2026-01-20")
        # are not silently converted into real date serial values, then
        # restore the default style so no formatting residue is left.
        $destRange.NumberFormat = "@"
        $destRange.Value = $value
        $destRange.Style = "Normal"
    } else {
        $destRange.Value = $value
    }
}

foreach ($c in $allCols) {
    $addr43 = $c + "43"
    $addr44 = $c + "44"

    $r43 = $ws.Range($addr43)
    $r44 = $ws.Range($addr44)

    $v43 = $r43.Value()
    $v44 = $r44.Value()

    # Skip columns where both rows already hold the same value (nothing to
    # swap, and re-writing it could leave needless formatting residue).
    if ($v43 -eq $v44) {
        continue
    }

    $isText = $textCols -contains $c

    Set-CellFrom $r43 $v44 $isText
    Set-CellFrom $r44 $v43 $isText
}
